$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F ("time_taken") - header styled like the other header cells
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# time_taken values for each data row (2-74)
$ws.Range("F2").Value = "2021-10-05 13:40:19.165861"
$ws.Range("F3").Value = "2021-10-05 13:40:19.165874"
$ws.Range("F4").Value = "2021-10-05 13:40:19.165879"
$ws.Range("F5").Value = "2021-10-05 13:40:19.165882"
$ws.Range("F6").Value = "2021-10-05 13:40:19.165885"
$ws.Range("F7").Value = "2021-10-05 13:40:19.165888"
$ws.Range("F8").Value = "2021-10-05 13:40:19.165891"
$ws.Range("F9").Value = "2021-10-05 13:40:19.165894"
$ws.Range("F10").Value = "2021-10-05 13:40:19.165898"
$ws.Range("F11").Value = "2021-10-05 13:40:19.165901"
$ws.Range("F12").Value = "2021-10-05 13:40:19.165903"
$ws.Range("F13").Value = "2021-10-05 13:40:19.165906"
$ws.Range("F14").Value = "2021-10-05 13:40:19.165909"
$ws.Range("F15").Value = "2021-10-05 13:40:19.165912"
$ws.Range("F16").Value = "2021-10-05 13:40:19.165915"
$ws.Range("F17").Value = "2021-10-05 13:40:19.165918"
$ws.Range("F18").Value = "2021-10-05 13:40:19.165921"
$ws.Range("F19").Value = "2021-10-05 13:40:19.165924"
$ws.Range("F20").Value = "2021-10-05 13:40:19.165927"
$ws.Range("F21").Value = "2021-10-05 13:40:19.165930"
$ws.Range("F22").Value = "2021-10-05 13:40:19.165933"
$ws.Range("F23").Value = "2021-10-05 13:40:19.165936"
$ws.Range("F24").Value = "2021-10-05 13:40:19.165939"
$ws.Range("F25").Value = "2021-10-05 13:40:19.165943"
$ws.Range("F26").Value = "2021-10-05 13:40:19.165946"
$ws.Range("F27").Value = "2021-10-05 13:40:19.165949"
$ws.Range("F28").Value = "2021-10-05 13:40:19.165952"
$ws.Range("F29").Value = "2021-10-05 13:40:19.165955"
$ws.Range("F30").Value = "2021-10-05 13:40:19.165958"
$ws.Range("F31").Value = "2021-10-05 13:40:19.165961"
$ws.Range("F32").Value = "2021-10-05 13:40:19.165964"
$ws.Range("F33").Value = "2021-10-05 13:40:19.165967"
$ws.Range("F34").Value = "2021-10-05 13:40:19.165970"
$ws.Range("F35").Value = "2021-10-05 13:40:19.165973"
$ws.Range("F36").Value = "2021-10-05 13:40:19.165976"
$ws.Range("F37").Value = "2021-10-05 13:40:19.165979"
$ws.Range("F38").Value = "2021-10-05 13:40:19.165982"
$ws.Range("F39").Value = "2021-10-05 13:40:19.165985"
$ws.Range("F40").Value = "2021-10-05 13:40:19.165988"
$ws.Range("F41").Value = "2021-10-05 13:40:19.165991"
$ws.Range("F42").Value = "2021-10-05 13:40:19.165994"
$ws.Range("F43").Value = "2021-10-05 13:40:19.165997"
$ws.Range("F44").Value = "2021-10-05 13:40:19.166000"
$ws.Range("F45").Value = "2021-10-05 13:40:19.166003"
$ws.Range("F46").Value = "2021-10-05 13:40:19.166006"
$ws.Range("F47").Value = "2021-10-05 13:40:19.166009"
$ws.Range("F48").Value = "2021-10-05 13:40:19.166012"
$ws.Range("F49").Value = "2021-10-05 13:40:19.166015"
$ws.Range("F50").Value = "2021-10-05 13:40:19.166018"
$ws.Range("F51").Value = "2021-10-05 13:40:19.166021"
$ws.Range("F52").Value = "2021-10-05 13:40:19.166024"
$ws.Range("F53").Value = "2021-10-05 13:40:19.166027"
$ws.Range("F54").Value = "2021-10-05 13:40:19.166030"
$ws.Range("F55").Value = "2021-10-05 13:40:19.166033"
$ws.Range("F56").Value = "2021-10-05 13:40:19.166036"
$ws.Range("F57").Value = "2021-10-05 13:40:19.166039"
$ws.Range("F58").Value = "2021-10-05 13:40:19.166042"
$ws.Range("F59").Value = "2021-10-05 13:40:19.166045"
$ws.Range("F60").Value = "2021-10-05 13:40:19.166048"
$ws.Range("F61").Value = "2021-10-05 13:40:19.166051"
$ws.Range("F62").Value = "2021-10-05 13:40:19.166054"
$ws.Range("F63").Value = "2021-10-05 13:40:19.166057"
$ws.Range("F64").Value = "2021-10-05 13:40:19.166060"
$ws.Range("F65").Value = "2021-10-05 13:40:19.166063"
$ws.Range("F66").Value = "2021-10-05 13:40:19.166067"
$ws.Range("F67").Value = "2021-10-05 13:40:19.166070"
$ws.Range("F68").Value = "2021-10-05 13:40:19.166073"
$ws.Range("F69").Value = "2021-10-05 13:40:19.166076"
$ws.Range("F70").Value = "2021-10-05 13:40:19.166079"
$ws.Range("F71").Value = "2021-10-05 13:40:19.166082"
$ws.Range("F72").Value = "2021-10-05 13:40:19.166085"
$ws.Range("F73").Value = "2021-10-05 13:40:19.166088"
$ws.Range("F74").Value = "2021-10-05 13:40:19.166091"
